$d = $word.ActiveDocument

# 1) CÓDIGO ÚNICO: 16 -> 17
$d.Content.Find.Execute("CÓDIGO ÚNICO: 16", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "CÓDIGO ÚNICO: 17", 2) | Out-Null

# 2) "  Destinado ao CRAS" -> "  pra fernanda"
$d.Content.Find.Execute("Destinado ao CRAS", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "pra fernanda", 2) | Out-Null

# 3) Remove the "Leite em pó" table row (last row of the only table)
$table = $d.Tables.Item(1)
for ($i = $table.Rows.Count; $i -ge 1; $i--) {
    $row = $table.Rows.Item($i)
    if ($row.Cells.Item(2).Range.Text -like "Leite em p*") {
        $row.Delete()
    }
}

# 4) Valor total R$ 16.88 -> R$ 13.78
$d.Content.Find.Execute("R`$ 16.88", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "R`$ 13.78", 2) | Out-Null

# 5) Update the timestamp at the bottom of the document
$d.Content.Find.Execute("Cortês/PE, 2024-10-14 18:40:17.205734+00:00.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Cortês/PE, 2024-10-16 19:46:17.837971+00:00.", 2) | Out-Null
